$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")

$ws_ALC.Range("H15").Value = 990.71606
$ws_ALC.Range("I15").Value = 990.71606
$ws_ALC.Range("K15").Value = 2972.14818
$ws_ALC.Range("M15").Value = -2803.14818

$ws_ALC.Range("H28").Value = 1873.2941
$ws_ALC.Range("J28").Value = 3935.7144
$ws_ALC.Range("L28").Value = 3935.7144
$ws_ALC.Range("N28").Value = -4905.7144

$ws_ALC.Range("H61").Value = 856.4
$ws_ALC.Range("J61").Value = 0
$ws_ALC.Range("L61").Value = 0
$ws_ALC.Range("N61").ClearContents()

$ws_ALC.Range("H114").Value = 17550
$ws_ALC.Range("J114").Value = 17550
$ws_ALC.Range("L114").Value = 17550
$ws_ALC.Range("N114").Value = -26228

$ws_ALC.Range("H137").Value = 100003300
$ws_ALC.Range("I137").Value = 66670508
$ws_ALC.Range("K137").Value = 200011524
$ws_ALC.Range("M137").Value = -200008974

$ws_ALC.Range("H138").Value = 3390.848
$ws_ALC.Range("J138").Value = 3839.879
$ws_ALC.Range("L138").Value = 11519.637
$ws_ALC.Range("N138").Value = -21799.637

$ws_ARM.Range("H2").Value = 3277.5833
$ws_ARM.Range("I2").Value = 3702
$ws_ARM.Range("K2").Value = 3702
$ws_ARM.Range("M2").Value = -3589

$ws_ARM.Range("H32").Value = 14937209
$ws_ARM.Range("I32").Value = 22227856
$ws_ARM.Range("K32").Value = 22227856
$ws_ARM.Range("M32").Value = -22227569

$ws_ARM.Range("H61").Value = 31253018
$ws_ARM.Range("I61").Value = 43480256
$ws_ARM.Range("J61").Value = 5632.4443
$ws_ARM.Range("K61").Value = 43480256
$ws_ARM.Range("L61").Value = 5632.4443
$ws_ARM.Range("M61").Value = -43480044
$ws_ARM.Range("N61").Value = -6056.4443

$ws_ARM.Range("H74").Value = 58891396
$ws_ARM.Range("I74").Value = 58891396
$ws_ARM.Range("J74").Value = 0
$ws_ARM.Range("K74").Value = 58891396
$ws_ARM.Range("L74").Value = 0
$ws_ARM.Range("M74").Value = -58890522
$ws_ARM.Range("N74").ClearContents()

$ws_ARM.Range("H77").Value = 58891396
$ws_ARM.Range("I77").Value = 58891396
$ws_ARM.Range("J77").Value = 0
$ws_ARM.Range("K77").Value = 294456980
$ws_ARM.Range("L77").Value = 0
$ws_ARM.Range("M77").Value = -294452612
$ws_ARM.Range("N77").ClearContents()

$ws_ARM.Range("H102").Value = 145527.64
$ws_ARM.Range("I102").Value = 183990.19
$ws_ARM.Range("J102").Value = 4498.3335
$ws_ARM.Range("K102").Value = 183990.19
$ws_ARM.Range("L102").Value = 4498.3335
$ws_ARM.Range("M102").Value = -182368.19
$ws_ARM.Range("N102").Value = -7742.3335

$ws_ARM.Range("H116").Value = 3277.5833
$ws_ARM.Range("I116").Value = 3702
$ws_ARM.Range("K116").Value = 3702
$ws_ARM.Range("M116").Value = -1408

$ws_ARM.Range("H132").Value = 26323344
$ws_ARM.Range("I132").Value = 8086.7715
$ws_ARM.Range("J132").Value = 333334660
$ws_ARM.Range("K132").Value = 24260.3145
$ws_ARM.Range("L132").Value = 1000003980
$ws_ARM.Range("M132").Value = -21730.3145
$ws_ARM.Range("N132").Value = -1000009040

$ws_ARM.Range("H136").Value = 31253018
$ws_ARM.Range("I136").Value = 43480256
$ws_ARM.Range("J136").Value = 5632.4443
$ws_ARM.Range("K136").Value = 130440768
$ws_ARM.Range("L136").Value = 16897.3329
$ws_ARM.Range("M136").Value = -130438218
$ws_ARM.Range("N136").Value = -21997.3329

$ws_BSM.Range("H3").Value = 3277.5833
$ws_BSM.Range("I3").Value = 3702
$ws_BSM.Range("K3").Value = 3702
$ws_BSM.Range("M3").Value = -3588

$ws_BSM.Range("H107").Value = 4501.4707
$ws_BSM.Range("I107").Value = 4342.2
$ws_BSM.Range("J107").Value = 4729
$ws_BSM.Range("K107").Value = 4342.2
$ws_BSM.Range("L107").Value = 4729
$ws_BSM.Range("M107").Value = -2422.2
$ws_BSM.Range("N107").Value = -8569

$ws_BSM.Range("H134").Value = 10497.5
$ws_BSM.Range("I134").Value = 8994.25
$ws_BSM.Range("J134").Value = 11249.125
$ws_BSM.Range("K134").Value = 26982.75
$ws_BSM.Range("L134").Value = 33747.375
$ws_BSM.Range("M134").Value = -24447.75
$ws_BSM.Range("N134").Value = -38817.375

$ws_CRP.Range("H16").Value = 1464.1875
$ws_CRP.Range("J16").Value = 1237.3334
$ws_CRP.Range("L16").Value = 1237.3334
$ws_CRP.Range("N16").Value = -1811.3334

$ws_CRP.Range("H22").Value = 12402.444
$ws_CRP.Range("I22").Value = 33520.332
$ws_CRP.Range("J22").Value = 1843.5
$ws_CRP.Range("K22").Value = 33520.332
$ws_CRP.Range("L22").Value = 1843.5
$ws_CRP.Range("M22").Value = -33170.332
$ws_CRP.Range("N22").Value = -2543.5

$ws_CRP.Range("H31").Value = 24394676
$ws_CRP.Range("I31").Value = 3416.2334
$ws_CRP.Range("K31").Value = 3416.2334
$ws_CRP.Range("M31").Value = -3121.2334

$ws_CRP.Range("H34").Value = 24394676
$ws_CRP.Range("I34").Value = 3416.2334
$ws_CRP.Range("K34").Value = 3416.2334
$ws_CRP.Range("M34").Value = -3214.2334

$ws_CRP.Range("H35").Value = 100007130
$ws_CRP.Range("I35").Value = 500005000
$ws_CRP.Range("J35").Value = 7663.375
$ws_CRP.Range("K35").Value = 500005000
$ws_CRP.Range("L35").Value = 7663.375
$ws_CRP.Range("M35").Value = -500004706
$ws_CRP.Range("N35").Value = -8251.375

$ws_CRP.Range("H86").Value = 3096.8125
$ws_CRP.Range("I86").Value = 2238.4
$ws_CRP.Range("J86").Value = 4527.5
$ws_CRP.Range("K86").Value = 2238.4
$ws_CRP.Range("L86").Value = 4527.5
$ws_CRP.Range("M86").Value = -1115.4
$ws_CRP.Range("N86").Value = -6773.5

$ws_CRP.Range("H89").Value = 3096.8125
$ws_CRP.Range("I89").Value = 2238.4
$ws_CRP.Range("J89").Value = 4527.5
$ws_CRP.Range("K89").Value = 11192
$ws_CRP.Range("L89").Value = 22637.5
$ws_CRP.Range("M89").Value = -5576
$ws_CRP.Range("N89").Value = -33869.5

$ws_CRP.Range("H107").Value = 1698.409
$ws_CRP.Range("I107").Value = 953.9091
$ws_CRP.Range("K107").Value = 953.9091
$ws_CRP.Range("M107").Value = 966.0909

$ws_CRP.Range("H113").Value = 1464.1875
$ws_CRP.Range("J113").Value = 1237.3334
$ws_CRP.Range("L113").Value = 1237.3334
$ws_CRP.Range("N113").Value = -5577.3334

$ws_GSM.Range("H43").Value = 16141.5
$ws_GSM.Range("J43").Value = 89999
$ws_GSM.Range("L43").Value = 89999
$ws_GSM.Range("N43").Value = -90301

$ws_GSM.Range("H46").Value = 30046
$ws_GSM.Range("I46").Value = 0
$ws_GSM.Range("J46").Value = 30046
$ws_GSM.Range("K46").Value = 0
$ws_GSM.Range("L46").Value = 30046
$ws_GSM.Range("M46").ClearContents()
$ws_GSM.Range("N46").Value = -30358

$ws_GSM.Range("H80").Value = 3978.4707
$ws_GSM.Range("I80").Value = 3869.4
$ws_GSM.Range("J80").Value = 4023.9167
$ws_GSM.Range("K80").Value = 3869.4
$ws_GSM.Range("L80").Value = 4023.9167
$ws_GSM.Range("M80").Value = -2871.4
$ws_GSM.Range("N80").Value = -6019.9167

$ws_GSM.Range("H83").Value = 3978.4707
$ws_GSM.Range("I83").Value = 3869.4
$ws_GSM.Range("J83").Value = 4023.9167
$ws_GSM.Range("K83").Value = 19347
$ws_GSM.Range("L83").Value = 20119.5835
$ws_GSM.Range("M83").Value = -14355
$ws_GSM.Range("N83").Value = -30103.5835

$ws_GSM.Range("H113").Value = 4117
$ws_GSM.Range("I113").Value = 3387
$ws_GSM.Range("K113").Value = 3387
$ws_GSM.Range("M113").Value = -1217

$ws_GSM.Range("H132").Value = 2285.5151
$ws_GSM.Range("I132").Value = 2106.9375
$ws_GSM.Range("J132").Value = 8000
$ws_GSM.Range("K132").Value = 6320.8125
$ws_GSM.Range("L132").Value = 24000
$ws_GSM.Range("M132").Value = -3790.8125
$ws_GSM.Range("N132").Value = -29060

$ws_LTW.Range("H40").Value = 3913.353
$ws_LTW.Range("I40").Value = 3502.6333
$ws_LTW.Range("J40").Value = 6993.75
$ws_LTW.Range("K40").Value = 3502.6333
$ws_LTW.Range("L40").Value = 6993.75
$ws_LTW.Range("M40").Value = -3366.6333
$ws_LTW.Range("N40").Value = -7265.75

$ws_LTW.Range("I46").Value = 993.75
$ws_LTW.Range("J46").Value = 2607.4
$ws_LTW.Range("K46").Value = 993.75
$ws_LTW.Range("L46").Value = 2607.4
$ws_LTW.Range("M46").Value = -805.75
$ws_LTW.Range("N46").Value = -2983.4

$ws_LTW.Range("H122").Value = 3600.9312
$ws_LTW.Range("I122").Value = 3137.08
$ws_LTW.Range("J122").Value = 6500
$ws_LTW.Range("K122").Value = 9411.24
$ws_LTW.Range("L122").Value = 19500
$ws_LTW.Range("M122").Value = -6961.24
$ws_LTW.Range("N122").Value = -24400

$ws_LTW.Range("H136").Value = 3946.5557
$ws_LTW.Range("I136").Value = 3946.5557
$ws_LTW.Range("K136").Value = 11839.6671
$ws_LTW.Range("M136").Value = -9289.667099999999

Write-Host "Applied all changes"